# Update "Take down" slide (slide 6): tidy up the run splits for a few
# lines in the content placeholder, then append two new slides at the
# end of the deck.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 6 ("Take down") - content placeholder text tidy-up
# ---------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$body6 = $slide6.Shapes.Item(2).TextFrame.TextRange

# Paragraph 3: "sudo iwconfig $card mode Managed"
#   merge the " $card " and "mode " runs into a single " $card mode " run
$para = $body6.Paragraphs(3)
$full = $para.Text
$idx = $full.IndexOf(" `$card ")
$mergeLen = (" `$card mode ").Length
$sub = $body6.Characters($para.Start + $idx, $mergeLen)
$sub.Text = " `$card mode "

# Paragraph 4: "sudo ifconfig $card up"
#   merge the " $" and "card " runs into a single " $card " run
$para = $body6.Paragraphs(4)
$full = $para.Text
$idx = $full.IndexOf(" `$")
$mergeLen = (" `$card ").Length
$sub = $body6.Characters($para.Start + $idx, $mergeLen)
$sub.Text = " `$card "

# Paragraph 5: "sudo service network-manager start"
#   split the "sudo " run into separate "sudo" and " " runs
$para = $body6.Paragraphs(5)
$sub = $body6.Characters($para.Start, 4)
$sub.Text = "sudo"

# ---------------------------------------------------------------------
# New slide 7: "Testing site: James Madison park"
# ---------------------------------------------------------------------
$slide7 = $p.Slides.Add($p.Slides.Count + 1, 2)
$slide7.Shapes.Item(1).TextFrame.TextRange.Text = "Testing site: James Madison park"

# ---------------------------------------------------------------------
# New slide 8: blank Title and Content slide
# ---------------------------------------------------------------------
$slide8 = $p.Slides.Add($p.Slides.Count + 1, 2)
